# Cinema/movies.xlsx - "Select format bug resolved"
#
# A new movie (CASA GUCCI) was added, a couple of ratings were corrected,
# and the whole table (plus some extra blank rows) was selected and
# reformatted as Text - which is why the numbers now come back in as
# shared strings instead of numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the used range plus extra blank rows below it and fix the
# "selecting a format" bug by forcing Text format on the whole block.
$rng = $ws.Range("A1:B30")
$rng.Select()
$rng.NumberFormat = "@"

# Re-enter the table contents (now stored as text because of the format
# applied above), including the newly added "CASA GUCCI" row and the
# couple of rating corrections. The order below matches the order the
# strings were (re-)typed in, which controls where they land in the
# shared-string table.
$ws.Range("A2").Value = "CASA GUCCI"
$ws.Range("B2").Value = "Rating: 6.9/10"

$ws.Range("A3").Value = "CLIFFORD, MARELE CAINE ROSU"
$ws.Range("B3").Value = "Rating: 6/10"

$ws.Range("A4").Value = "ENCANTO"
$ws.Range("B4").Value = "Rating: 4.7/10"

$ws.Range("A5").Value = "HAI SA CANTAM DIN NOU!"
$ws.Range("B5").Value = "Rating: 7.6/10"

$ws.Range("A6").Value = "MATRIX RENASTEREA"
$ws.Range("B6").Value = "Rating: 5.7/10"

$ws.Range("A7").Value = "OMUL-PĂIANJEN: NICIUN DRUM SPRE CASĂ"
$ws.Range("B7").Value = "Rating: 8.8/10"

$ws.Range("A8").Value = "SPENCER"
$ws.Range("B8").Value = "Rating: 6.8/10"

$ws.Range("A9").Value = "TABARA"
$ws.Range("B9").Value = "Rating: 2.5/10"

$ws.Range("A10").Value = "TE URASC, TE IUBESC"
$ws.Range("B10").Value = "Rating: 6.4/10"

$ws.Range("A11").Value = "THE KING'S MAN: ÎNCEPUTUL"
$ws.Range("B11").Value = "Rating: 6.8/10"

$ws.Range("A12").Value = "VANATORII DE FANTOME: MOSTENIREA"
$ws.Range("B12").Value = "Rating: 7.4/10"

$ws.Range("A13").Value = "THE KING'S MAN: ÎNCEPUTUL"

$ws.Range("A14").Value = "VANATORII DE FANTOME: MOSTENIREA"

$ws.Range("A1").Value = "355"
$ws.Range("B1").Value = "Rating: 4.8/10"

# Auto-fit the two columns to their new content (target widths: ~40 and
# ~13.14 characters, as a real AutoFit would compute for this data/font).
$ws.Columns("A").ColumnWidth = 39.166666666666664
$ws.Columns("B").ColumnWidth = 12.307291666666666

# Final click that ends up as the saved selection.
$ws.Range("D7").Select()
